# Update the "Main.xlsx" Rules worksheet: cell E8 changes from "Good Morning"
# to "GIT UPDATE", and the cell becomes the active/selected cell, matching
# the sheetView's new <selection activeCell="E8" sqref="E8"/>.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"

$ws.Activate()
$ws.Range("E8").Select()
